# "edit path and edit test case"
# - Column C test-data JSON strings get a trailing "}" (closing brace) added.
# - Column D expected-result strings are translated from English to Chinese.
# - Active-sheet selection moves from D7 to D5.
# - Column D is widened.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login_data")   # "login_data" is the active/selected sheet

# --- Fix the JSON-looking request bodies in column C (close the '{' with a '}') ---
$ws.Range("C2").Value = "{'email':'ktp0215926300','password':'test123','remember':0}"
$ws.Range("C3").Value = "{'email':'','password':'test123','remember':0}"
$ws.Range("C4").Value = "{'email':'ktp0215926300','password':'','remember':0}"
$ws.Range("C5").Value = "{'email':'ktp0215926300','password':'error123','remember':0}"

# --- Translate the expected-result messages in column D ---
$ws.Range("D3").Value = "用户名不能为空"
$ws.Range("D4").Value = "密码不能为空"
$ws.Range("D5").Value = "密码错误, 你还可以尝试4次"

# --- Widen column D to fit the new, longer Chinese text ---
$ws.Columns.Item(4).ColumnWidth = 29.15

# --- Move the active selection to D5 ---
[void]$ws.Range("D5").Select()
